$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list snapshot: refresh Price (D) and Volume(1h) (E) columns
# with new scraped values. Cells hold literal text (e.g. trailing zeros,
# "%" suffix) so each target cell is first formatted as Text to stop Excel
# from auto-converting the typed string into a Number/Percentage.
$updates = @{
    "D2" = "302.16"
    "E2" = "2.43%"
    "D3" = "35.19"
    "E3" = "13.52%"
    "D4" = "5.153"
    "E4" = "4.97%"
    "D5" = "0.07769"
    "E5" = "4.38%"
    "D6" = "2.355"
    "E6" = "8.20%"
    "D7" = "8.036"
    "D8" = "3.950"
    "E8" = "5.32%"
    "D9" = "0.9316"
    "E9" = "1.95%"
    "D10" = "0.09956"
    "E10" = "12.20%"
    "D11" = "0.1791"
    "E11" = "4.59%"
    "D12" = "0.08603"
    "E12" = "3.66%"
    "D13" = "0.03312"
    "E13" = "5.27%"
    "D14" = "0.09918"
    "E14" = "-1.44%"
    "D15" = "0.001487"
    "E15" = "-1.67%"
    "D16" = "0.005772"
    "E16" = "0.97%"
    "D18" = "2.146"
    "E18" = "3.37%"
    "D20" = "0.1303"
    "E20" = "0.78%"
    "D21" = "4.298"
    "E21" = "8.29%"
    "D22" = "0.2302"
    "E22" = "9.53%"
    "D23" = "0.04534"
    "E23" = "-0.49%"
    "E24" = "0.16%"
    "D25" = "0.004378"
    "E25" = "-5.41%"
    "E26" = "-0.01%"
    "E27" = "-0.01%"
    "D39" = "0.01780"
    "E39" = "10.13%"
    "D40" = "0.04797"
    "E40" = "6.84%"
    "D41" = "0.007779"
    "E41" = "6.78%"
    "D42" = "0.1412"
    "E42" = "5.98%"
    "D43" = "0.006830"
    "E43" = "-23.83%"
    "D44" = "0.002073"
    "E44" = "5.52%"
    "D45" = "0.009443"
    "E45" = "3.27%"
    "D46" = "0.00006113"
    "E46" = "-0.02%"
    "E47" = "0.00%"
    "D48" = "2.985"
    "E48" = "33.51%"
    "E49" = "0.08%"
    "E50" = "0.00%"
    "E51" = "0.00%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}

